$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titleText = "Play Battleship Direct Hit Free - Exciting Online Slot Game"

# ---------------------------------------------------------------------------
# 1) Swap the final (italic) paragraph's copy text for the new image-prompt
#    text, keeping its run formatting (italic) intact. InsertXML (rather
#    than Find/Replace) is used so straight quotes/apostrophes in the new
#    copy are not silently "smart-quoted".
# ---------------------------------------------------------------------------
$newCopy = 'Please create an image featuring a happy Maya warrior wearing glasses for the game "Battleship Direct Hit". The image should be in a cartoon style and should capture the essence of the game''s naval battle theme in a fun and engaging way. It should also feature the game''s title prominently. Be creative and use vibrant colors and dynamic imagery to attract players to this exciting slot game.'

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastRange = $lastPara.Range
$newLastXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$newCopy</w:t></w:r></w:p>"
$lastRange.InsertXML($newLastXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Battleship Direct Hit Free - Exciting
#    Online Slot Game" paragraph that used to sit right before the final
#    italic paragraph (now immediately before it).
# ---------------------------------------------------------------------------
$dupParaIndex = 0
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText) {
        $dupParaIndex = $i
        break
    }
}

if ($dupParaIndex -gt 0) {
    $d.Paragraphs.Item($dupParaIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Battleship Direct Hit Free -
#    Exciting Online Slot Game").
# ---------------------------------------------------------------------------
$titleParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText) {
        $titleParaIndex = $i
        break
    }
}

$titlePara = $d.Paragraphs.Item($titleParaIndex)
# Position right before the paragraph's own end-of-paragraph mark so the
# inserted XML becomes its own, brand new paragraph (no inherited pStyle).
$insertPos = $titlePara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$metaXml = "<w:p $wNs>" +
           "<w:r/>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
           "<w:r><w:t>: Read our expert review of Battleship Direct Hit and play for free. Impressive graphics, Megaways, and free spins make for an exciting gaming experience.</w:t></w:r>" +
           "</w:p>"
$insertRange.InsertXML($metaXml)
